$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1697.5
$ws.Range("J17").Value = 2096.6667
$ws.Range("L17").Value = 6290.000100000001
$ws.Range("N17").Value = -6626.000100000001

$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").Value = ""

$ws.Range("H43").Value = 7500
$ws.Range("J43").Value = 7500
$ws.Range("L43").Value = 7500
$ws.Range("N43").Value = -7638

$ws.Range("H138").Value = 3506.4614
$ws.Range("J138").Value = 3642.889
$ws.Range("L138").Value = 10928.667
$ws.Range("N138").Value = -21208.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""

$ws.Range("H76").Value = 17821.75
$ws.Range("J76").Value = 17821.75
$ws.Range("L76").Value = 17821.75
$ws.Range("N76").Value = -18497.75

$ws.Range("H79").Value = 17821.75
$ws.Range("J79").Value = 17821.75
$ws.Range("L79").Value = 17821.75
$ws.Range("N79").Value = -20161.75

$ws.Range("H132").Value = 1440.7142
$ws.Range("I132").Value = 1514.1666
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 4542.4998
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -2012.4998
$ws.Range("N132").Value = -8060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3104.75
$ws.Range("I20").Value = 2998
$ws.Range("K20").Value = 2998
$ws.Range("M20").Value = -2751

$ws.Range("H95").Value = 16314.6
$ws.Range("J95").Value = 16314.6
$ws.Range("L95").Value = 16314.6
$ws.Range("N95").Value = -21806.6

$ws.Range("H97").Value = 19999
$ws.Range("I97").Value = 19999
$ws.Range("K97").Value = 19999
$ws.Range("M97").Value = -19008

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = ""
$ws.Range("N134").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 2016.909
$ws.Range("I32").Value = 538.2
$ws.Range("J32").Value = 3249.1667
$ws.Range("K32").Value = 538.2
$ws.Range("L32").Value = 3249.1667
$ws.Range("M32").Value = -222.2
$ws.Range("N32").Value = -3881.1667

$ws.Range("H35").Value = 1774.5
$ws.Range("J35").Value = 2997.5
$ws.Range("L35").Value = 2997.5
$ws.Range("N35").Value = -3585.5

$ws.Range("H45").Value = 24749
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 24749
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 24749
$ws.Range("M45").Value = ""
$ws.Range("N45").Value = -25935

$ws.Range("H93").Value = 45302.332
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = ""

$ws.Range("H96").Value = 19996.334
$ws.Range("J96").Value = 19996.334
$ws.Range("L96").Value = 19996.334
$ws.Range("N96").Value = -25488.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 269.33334
$ws.Range("J107").Value = 269.33334
$ws.Range("L107").Value = 808.0000200000001
$ws.Range("N107").Value = -4648.00002

$ws.Range("H134").Value = 2810.8333
$ws.Range("I134").Value = 1458
$ws.Range("K134").Value = 4374
$ws.Range("M134").Value = 696

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = ""

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7900
$ws.Range("I46").Value = 7900
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 7900
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -7712
$ws.Range("N46").Value = ""

$ws.Range("H68").Value = 2400
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""

$ws.Range("H71").Value = 2400
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""

$ws.Range("H76").Value = 32000
$ws.Range("J76").Value = 32000
$ws.Range("L76").Value = 32000
$ws.Range("N76").Value = -32676

$ws.Range("H79").Value = 32000
$ws.Range("J79").Value = 32000
$ws.Range("L79").Value = 32000
$ws.Range("N79").Value = -34340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 500
$ws.Range("J29").Value = 500
$ws.Range("L29").Value = 500
$ws.Range("N29").Value = -1080

$ws.Range("H62").Value = 200000
$ws.Range("J62").Value = 200000
$ws.Range("L62").Value = 200000
$ws.Range("N62").Value = -201248

$ws.Range("H65").Value = 200000
$ws.Range("J65").Value = 200000
$ws.Range("L65").Value = 1000000
$ws.Range("N65").Value = -1006240

$ws.Range("H69").Value = 16199.4
$ws.Range("J69").Value = 16199.4
$ws.Range("L69").Value = 16199.4
$ws.Range("N69").Value = -17697.4

$ws.Range("H72").Value = 16199.4
$ws.Range("J72").Value = 16199.4
$ws.Range("L72").Value = 48598.2
$ws.Range("N72").Value = -56086.2

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").Value = ""

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").Value = ""

$ws.Range("H100").Value = 1916.3334
$ws.Range("I100").Value = 1000
$ws.Range("J100").Value = 2374.5
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 4749
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -5831

$ws.Range("H104").Value = 26158.334
$ws.Range("J104").Value = 27390
$ws.Range("L104").Value = 27390
$ws.Range("N104").Value = -34378

$ws.Range("H126").Value = 9333
$ws.Range("I126").Value = 9333
$ws.Range("K126").Value = 27999
$ws.Range("M126").Value = -25529
